$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Summary of the change (per the commit "updated jax fc code"):
#
# Item YT427896-805 block (rows 6-13):
#   rows 6-9   : Date '11/19/2020' Qty 4  ->  Date '1/18/2021'  Qty 3
#   rows 10-13 : Date '12/10/2020' Qty 3  ->  Date '11/19/2020' Qty 4
#
# Item YT427897-805 block (rows 22-29):
#   rows 22-25 : Date '11/9/2020'  Qty 4  ->  Date '1/18/2021'  Qty 3
#   rows 26-29 : Date '12/10/2020' Qty 3  ->  Date '11/9/2020'  Qty 4
#
# Extended Price (column F) = Price (column E, constant 9119) * Qty, so
# it moves between 36476 (qty 4) and 27357 (qty 3) accordingly.
#
# The Date column stores plain text (shared strings), not real Excel
# dates. We must avoid letting Excel auto-convert the text into a date
# serial number. The safest way to do this without disturbing styles is
# to copy the date value from an existing cell that already holds the
# desired text (same-size range-to-range copy), rather than typing a
# literal date string into the Value property.
# -------------------------------------------------------------------------

# Step 1: Before anything else, propagate the current (pre-edit) date
# text of rows 6-9 ('11/19/2020') and rows 22-25 ('11/9/2020') down into
# rows 10-13 and 26-29 respectively, since those source cells are about
# to be overwritten with '1/18/2021' in Step 2.
$ws.Range("B6:B9").Copy($ws.Range("B10:B13"))
$ws.Range("B22:B25").Copy($ws.Range("B26:B29"))

# Step 2: Now overwrite rows 6-9 and 22-25 with '1/18/2021', copied from
# the existing cells B18:B21 which already contain that exact text as a
# shared string (this avoids Excel re-interpreting the literal string
# "1/18/2021" as a date serial number and avoids adding new styles).
$ws.Range("B18:B21").Copy($ws.Range("B6:B9"))
$ws.Range("B18:B21").Copy($ws.Range("B22:B25"))

# Step 3: Update the Qty column (C) to match the new dates.
$ws.Range("C6:C9").Value = 3
$ws.Range("C22:C25").Value = 3
$ws.Range("C10:C13").Value = 4
$ws.Range("C26:C29").Value = 4

# Step 4: Update the Extended Price column (F) = Price * Qty (Price is a
# constant 9119 for these rows).
$ws.Range("F6:F9").Value = 27357
$ws.Range("F22:F25").Value = 27357
$ws.Range("F10:F13").Value = 36476
$ws.Range("F26:F29").Value = 36476
